$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a pure number (e.g. "1.01") need the
# number format forced to Text first, otherwise Excel auto-converts the
# assigned string into a numeric value instead of keeping it as text.
$textCells = @(
    "D4",
    "D5",
    "D9",
    "D10",
    "D14",
    "D15",
    "D16",
    "D19",
    "D21",
    "D22",
    "D23",
    "D25",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D38",
    "D39",
    "D41",
    "D43",
    "D45",
    "D46",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.645.58"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.641.84"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "214.95"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "19.05"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.871.89"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.649.53"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "65.03"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "26.692.85"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "215.38"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "9.48"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  +12.47%  "
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").Value = "15.69"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").Value = "1.279.87"
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  +6.63%  "
$ws.Range("D39").Value = "0.830"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").Value = "1.781.93"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "91.24"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "59.90"
$ws.Range("E46").Value = "  +8.79%  "
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "7.82"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  -0.34%  "

# Restore the default "Normal" cell style on the cells we forced to Text
# so no stray number-format styling is left behind on the cell.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
